# Auto-generated: apply updated market/profit figures to the per-job "Leve" sheets.
# Values come from a scheduled data-refresh run (see commit message); only the
# H:N price/profit columns change, row identities (A:G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 174.28572
$ws.Range("I5").Value = 53.333332
$ws.Range("J5").Value = 265
$ws.Range("K5").Value = 53.333332
$ws.Range("L5").Value = 265
$ws.Range("M5").Value = 61.666668
$ws.Range("N5").Value = -495
$ws.Range("H9").Value = 294.25
$ws.Range("I9").Value = 186
$ws.Range("K9").Value = 186
$ws.Range("M9").Value = -17
$ws.Range("H11").Value = 150.6
$ws.Range("I11").Value = 150.6
$ws.Range("K11").Value = 150.6
$ws.Range("M11").Value = -10.59999999999999
$ws.Range("H18").Value = 3197.25
$ws.Range("I18").Value = 3525.4285
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 3525.4285
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = -3241.4285
$ws.Range("N18").Value = -1468
$ws.Range("H64").Value = 10012.5
$ws.Range("I64").Value = 4350
$ws.Range("J64").Value = 11900
$ws.Range("K64").Value = 4350
$ws.Range("L64").Value = 11900
$ws.Range("M64").Value = -4102
$ws.Range("N64").Value = -12396
$ws.Range("H67").Value = 10012.5
$ws.Range("I67").Value = 4350
$ws.Range("J67").Value = 11900
$ws.Range("K67").Value = 4350
$ws.Range("L67").Value = 11900
$ws.Range("M67").Value = -3492
$ws.Range("N67").Value = -13616
$ws.Range("H92").Value = 627
$ws.Range("I92").Value = 627
$ws.Range("K92").Value = 627
$ws.Range("M92").Value = 621
$ws.Range("H98").Value = 3674.4614
$ws.Range("I98").Value = 3914.0833
$ws.Range("K98").Value = 3914.0833
$ws.Range("M98").Value = -2416.0833
$ws.Range("H122").Value = 3674.4614
$ws.Range("I122").Value = 3914.0833
$ws.Range("K122").Value = 11742.2499
$ws.Range("M122").Value = -9292.249899999999
$ws.Range("H127").Value = 916.7778
$ws.Range("I127").Value = 843.875
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 2531.625
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 2428.375
$ws.Range("N127").Value = -14420
$ws.Range("H132").Value = 1383.5714
$ws.Range("I132").Value = 1402.35
$ws.Range("J132").Value = 1008
$ws.Range("K132").Value = 4207.049999999999
$ws.Range("L132").Value = 3024
$ws.Range("M132").Value = -1677.049999999999
$ws.Range("N132").Value = -8084
$ws.Range("H138").Value = 30306334
$ws.Range("J138").Value = 43481692
$ws.Range("L138").Value = 130445076
$ws.Range("N138").Value = -130455356

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3691.1
$ws.Range("I2").Value = 2864.125
$ws.Range("J2").Value = 6999
$ws.Range("K2").Value = 2864.125
$ws.Range("L2").Value = 6999
$ws.Range("M2").Value = -2751.125
$ws.Range("N2").Value = -7225
$ws.Range("H61").Value = 38471000
$ws.Range("I61").Value = 50007800
$ws.Range("K61").Value = 50007800
$ws.Range("M61").Value = -50007588
$ws.Range("H116").Value = 3691.1
$ws.Range("I116").Value = 2864.125
$ws.Range("J116").Value = 6999
$ws.Range("K116").Value = 2864.125
$ws.Range("L116").Value = 6999
$ws.Range("M116").Value = -570.125
$ws.Range("N116").Value = -11587
$ws.Range("H132").Value = 4888.05
$ws.Range("I132").Value = 4192.3057
$ws.Range("J132").Value = 11149.75
$ws.Range("K132").Value = 12576.9171
$ws.Range("L132").Value = 33449.25
$ws.Range("M132").Value = -10046.9171
$ws.Range("N132").Value = -38509.25
$ws.Range("H136").Value = 38471000
$ws.Range("I136").Value = 50007800
$ws.Range("K136").Value = 150023400
$ws.Range("M136").Value = -150020850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3691.1
$ws.Range("I3").Value = 2864.125
$ws.Range("J3").Value = 6999
$ws.Range("K3").Value = 2864.125
$ws.Range("L3").Value = 6999
$ws.Range("M3").Value = -2750.125
$ws.Range("N3").Value = -7227
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H134").Value = 2299.814
$ws.Range("I134").Value = 2326.4614
$ws.Range("J134").Value = 2040
$ws.Range("K134").Value = 6979.3842
$ws.Range("L134").Value = 6120
$ws.Range("M134").Value = -4444.3842
$ws.Range("N134").Value = -11190

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2559.6191
$ws.Range("I16").Value = 2018.8572
$ws.Range("J16").Value = 3641.1428
$ws.Range("K16").Value = 2018.8572
$ws.Range("L16").Value = 3641.1428
$ws.Range("M16").Value = -1731.8572
$ws.Range("N16").Value = -4215.1428
$ws.Range("H28").Value = 49984.5
$ws.Range("J28").Value = 49984.5
$ws.Range("L28").Value = 49984.5
$ws.Range("N28").Value = -50474.5
$ws.Range("H31").Value = 7021.067
$ws.Range("J31").Value = 7477
$ws.Range("L31").Value = 7477
$ws.Range("N31").Value = -8067
$ws.Range("H34").Value = 7021.067
$ws.Range("J34").Value = 7477
$ws.Range("L34").Value = 7477
$ws.Range("N34").Value = -7881
$ws.Range("H58").Value = 6397.6665
$ws.Range("I58").Value = 2096.1667
$ws.Range("J58").Value = 9838.866
$ws.Range("K58").Value = 2096.1667
$ws.Range("L58").Value = 9838.866
$ws.Range("M58").Value = -1893.1667
$ws.Range("N58").Value = -10244.866
$ws.Range("H99").Value = 4800
$ws.Range("I99").Value = 4800
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4800
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3302
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 2559.6191
$ws.Range("I113").Value = 2018.8572
$ws.Range("J113").Value = 3641.1428
$ws.Range("K113").Value = 2018.8572
$ws.Range("L113").Value = 3641.1428
$ws.Range("M113").Value = 151.1428000000001
$ws.Range("N113").Value = -7981.1428
$ws.Range("H126").Value = 4800
$ws.Range("I126").Value = 4800
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14400
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11930
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 6397.6665
$ws.Range("I136").Value = 2096.1667
$ws.Range("J136").Value = 9838.866
$ws.Range("K136").Value = 6288.500100000001
$ws.Range("L136").Value = 29516.598
$ws.Range("M136").Value = -3738.500100000001
$ws.Range("N136").Value = -34616.598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1228
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 1200
$ws.Range("M33").Value = -917
$ws.Range("H44").Value = 13000002
$ws.Range("I44").Value = 13000002
$ws.Range("K44").Value = 39000006
$ws.Range("M44").Value = -38999608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4333
$ws.Range("I4").Value = 3999.5
$ws.Range("K4").Value = 3999.5
$ws.Range("M4").Value = -3887.5
$ws.Range("H5").Value = 999
$ws.Range("I5").Value = 999
$ws.Range("K5").Value = 999
$ws.Range("M5").Value = -887
$ws.Range("H70").Value = 10499.333
$ws.Range("I70").Value = 8600
$ws.Range("K70").Value = 8600
$ws.Range("M70").Value = -8330
$ws.Range("H73").Value = 10499.333
$ws.Range("I73").Value = 8600
$ws.Range("K73").Value = 8600
$ws.Range("M73").Value = -7664
$ws.Range("H132").Value = 3628.551
$ws.Range("I132").Value = 3086.4856
$ws.Range("K132").Value = 9259.4568
$ws.Range("M132").Value = -6729.4568
$ws.Range("H133").Value = 83639.75
$ws.Range("J133").Value = 83639.75
$ws.Range("L133").Value = 83639.75
$ws.Range("N133").Value = -93759.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1432.25
$ws.Range("I22").Value = 969.7143
$ws.Range("K22").Value = 969.7143
$ws.Range("M22").Value = -674.7143
$ws.Range("H27").Value = 1432.25
$ws.Range("I27").Value = 969.7143
$ws.Range("K27").Value = 969.7143
$ws.Range("M27").Value = -862.7143
$ws.Range("H68").Value = 14139.286
$ws.Range("I68").Value = 14196
$ws.Range("K68").Value = 14196
$ws.Range("M68").Value = -13447
$ws.Range("H71").Value = 14139.286
$ws.Range("I71").Value = 14196
$ws.Range("K71").Value = 70980
$ws.Range("M71").Value = -67236
$ws.Range("H82").Value = 6590.3105
$ws.Range("I82").Value = 6884.647
$ws.Range("J82").Value = 6173.3335
$ws.Range("K82").Value = 6884.647
$ws.Range("L82").Value = 6173.3335
$ws.Range("M82").Value = -6523.647
$ws.Range("N82").Value = -6895.3335
$ws.Range("H85").Value = 6590.3105
$ws.Range("I85").Value = 6884.647
$ws.Range("J85").Value = 6173.3335
$ws.Range("K85").Value = 6884.647
$ws.Range("L85").Value = 6173.3335
$ws.Range("M85").Value = -5636.647
$ws.Range("N85").Value = -8669.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 5157
$ws.Range("J62").Value = 5705
$ws.Range("L62").Value = 5705
$ws.Range("N62").Value = -6953
$ws.Range("H65").Value = 5157
$ws.Range("J65").Value = 5705
$ws.Range("L65").Value = 28525
$ws.Range("N65").Value = -34765
$ws.Range("H96").Value = 3575.5
$ws.Range("J96").Value = 3712.5
$ws.Range("L96").Value = 3712.5
$ws.Range("N96").Value = -6458.5
$ws.Range("H122").Value = 4238.4814
$ws.Range("I122").Value = 2963.4666
$ws.Range("J122").Value = 5832.25
$ws.Range("K122").Value = 8890.399800000001
$ws.Range("L122").Value = 17496.75
$ws.Range("M122").Value = -6440.399800000001
$ws.Range("N122").Value = -22396.75

